$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trivia questions for the "General Knowledge" column (D), rows 2-11.
$questions = @(
    "What is Joey Chesnutt's record for most hot dogs eaten in a ten minute period?",
    "What U.S. city has the most breweries per capita?",
    "What is the most widely produced crop in the world?",
    "How many Goldfish Crackers are in one serving?",
    "What ingredient gives Malört",
    "What is the scale used to measure the spice level of peppers?",
    "What is the oldest brewery in the United States?",
    "How many gallons of water does it take to make a 1/3 pound hamburger?",
    "What are the five flavors of an original Rainbow Cone?",
    "How many ingredients are in a Manhattan (including the garnish)?"
)

# Row heights that Excel computed after the longer, wrapped question text was
# entered (row 6 and row 10 already had enough height from column E content).
$rowHeights = @{
    2  = 85
    3  = 51
    4  = 51
    5  = 51
    7  = 68
    8  = 51
    9  = 68
    11 = 68
}

for ($i = 0; $i -lt $questions.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $questions[$i]
    if ($rowHeights.ContainsKey($row)) {
        $ws.Rows.Item($row).RowHeight = $rowHeights[$row]
    }
}

# Scroll the view down a couple of rows and move the selection, matching the
# saved workbook view state.
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$ws.Range("D12").Select()
